$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells to preserve as text (matching original inlineStr type)
$textCells = @("D5", "D6", "D9", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D24", "D25", "D26", "D28", "D29", "D32", "D35", "D36", "D37", "D41", "D43", "D44", "D46", "D47", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range("D2").Value = "34.812.63"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "1.816.92"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "230.03"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -7.77%  "
$ws.Range("D9").Value = "0.322"
$ws.Range("E9").Value = "  +3.60%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "2.079.34"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.855.60"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "11.23"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "0.666"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "4.59"
$ws.Range("E16").Value = "  -2.53%  "
$ws.Range("D17").Value = "34.799.01"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "69.59"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  -1.64%  "
$ws.Range("D20").Value = "240.65"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "11.97"
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "4.66"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "2.26"
$ws.Range("E24").Value = "  +0.70%  "
$ws.Range("D25").Value = "173.08"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "7.74"
$ws.Range("E26").Value = "  -1.50%  "
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").Value = "17.27"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("D29").Value = "1.51"
$ws.Range("E29").Value = "  -5.15%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "0.0546"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("E34").Value = "  +10.45%  "
$ws.Range("D35").Value = "1.83"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "0.688"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").Value = "91.45"
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("E38").Value = "  +6.32%  "
$ws.Range("D39").Value = "1.339.08"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("D41").Value = "0.970"
$ws.Range("E41").Value = "  -2.97%  "
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("D43").Value = "2.42"
$ws.Range("E43").Value = "  -1.53%  "
$ws.Range("D44").Value = "14.19"
$ws.Range("E44").Value = "  -6.34%  "
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("D46").Value = "0.0523"
$ws.Range("E46").Value = "  +2.95%  "
$ws.Range("D47").Value = "6.17"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").Value = "1.995.11"
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "49.20"
$ws.Range("E51").Value = "  -0.30%  "

# Reset style back to Normal (remove temporary NumberFormat) for the forced-text cells
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}